# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.860.19"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$ws.Range("D3").Value = "2.047.65"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'245.02"
$ws.Range("E5").Value = "  -1.32%  "

# Row 6
$ws.Range("D6").Value = "'0.653"
$ws.Range("E6").Value = "  -1.71%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.12"
$ws.Range("E8").Value = "  -3.32%  "

# Row 9
$ws.Range("D9").Value = "'58.89"
$ws.Range("E9").Value = "  -1.69%  "

# Row 10
$ws.Range("E10").Value = "  -3.72%  "

# Row 11
$ws.Range("D11").Value = "'0.0774"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("E12").Value = "  +1.90%  "

# Row 13
$ws.Range("D13").Value = "'15.04"
$ws.Range("E13").Value = "  -4.49%  "

# Row 14
$ws.Range("E14").Value = "  +5.07%  "

# Row 15
$ws.Range("D15").Value = "2.347.39"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("D16").Value = "'5.55"
$ws.Range("E16").Value = "  -3.23%  "

# Row 17
$ws.Range("D17").Value = "2.080.38"
$ws.Range("E17").Value = "  +1.69%  "

# Row 18
$ws.Range("D18").Value = "36.845.13"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").Value = "'17.46"
$ws.Range("E19").Value = "  -1.86%  "

# Row 20
$ws.Range("D20").Value = "'73.03"
$ws.Range("E20").Value = "  -2.65%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22
$ws.Range("D22").Value = "'5.39"
$ws.Range("E22").Value = "  +1.17%  "

# Row 23
$ws.Range("D23").Value = "'235.66"
$ws.Range("E23").Value = "  -0.63%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  +7.04%  "

# Row 27
$ws.Range("E27").Value = "  +1.60%  "

# Row 28
$ws.Range("D28").Value = "'168.61"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").Value = "'20.09"
$ws.Range("E29").Value = "  +0.50%  "

# Row 30
$ws.Range("E30").Value = "  +14.33%  "

# Row 31
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("E32").Value = "  -0.25%  "

# Row 33
$ws.Range("D33").Value = "'4.78"
$ws.Range("E33").Value = "  +6.20%  "

# Row 34
$ws.Range("D34").Value = "'0.0613"
$ws.Range("E34").Value = "  -1.63%  "

# Row 35
$ws.Range("D35").Value = "'2.34"
$ws.Range("E35").Value = "  +6.05%  "

# Row 36
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("E37").Value = "  +5.65%  "

# Row 38
$ws.Range("D38").Value = "'0.0849"
$ws.Range("E38").Value = "  -5.42%  "

# Row 39
$ws.Range("D39").Value = "'1.31"
$ws.Range("E39").Value = "  -2.18%  "

# Row 40
$ws.Range("E40").Value = "  +0.40%  "

# Row 41
$ws.Range("E41").Value = "  -6.93%  "

# Row 42
$ws.Range("D42").Value = "'4.88"
$ws.Range("E42").Value = "  -4.22%  "

# Row 43
$ws.Range("E43").Value = "  +0.95%  "

# Row 44
$ws.Range("E44").Value = "  -10.57%  "

# Row 45
$ws.Range("D45").Value = "'96.59"
$ws.Range("E45").Value = "  +0.66%  "

# Row 46
$ws.Range("D46").Value = "'16.67"
$ws.Range("E46").Value = "  -4.08%  "

# Row 47
$ws.Range("D47").Value = "1.301.79"
$ws.Range("E47").Value = "  +1.45%  "

# Row 48
$ws.Range("E48").Value = "  -4.31%  "

# Row 49
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  -2.15%  "

# Row 50
$ws.Range("D50").Value = "'6.74"
$ws.Range("E50").Value = "  -0.45%  "

# Row 51
$ws.Range("D51").Value = "2.233.91"
$ws.Range("E51").Value = "  +0.06%  "

